# Atualizacao de bases das ligas, do dia: 2024-01-29 as 23-06
#
# The source data rows got re-sorted/re-paired; for several adjacent row
# pairs the entire match record (every column except the running index in
# column A) needs to be swapped between the two rows.
#
# Row pairs (1-based worksheet rows) whose B:AC data must be exchanged:
#   474 <-> 475
#   483 <-> 484
#   506 <-> 507
#   540 <-> 541
#   547 <-> 548
#   553 <-> 554
#   571 <-> 572

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(474, 475),
    @(483, 484),
    @(506, 507),
    @(540, 541),
    @(547, 548),
    @(553, 554),
    @(571, 572)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
